$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = 'Volume 30   Number  18'
$ws.Cells.Item(9, 3).Value = 'Report Covering the Week  5/1/2023  Through  5/7/2023'
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = -100
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = -60
$ws.Cells.Item(14, 14).Value = -88.888888888888
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 13
$ws.Cells.Item(15, 10).Value = 7
$ws.Cells.Item(15, 11).Value = 85.714285714285
$ws.Cells.Item(15, 12).Value = -31.578947368421
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = -62.857142857142
$ws.Cells.Item(16, 3).Value = 9
$ws.Cells.Item(16, 4).Value = 4
$ws.Cells.Item(16, 5).Value = 125
$ws.Cells.Item(16, 6).Value = 19
$ws.Cells.Item(16, 7).Value = 28
$ws.Cells.Item(16, 8).Value = -32.142857142857
$ws.Cells.Item(16, 9).Value = 87
$ws.Cells.Item(16, 10).Value = 105
$ws.Cells.Item(16, 11).Value = -17.142857142857
$ws.Cells.Item(16, 12).Value = 29.850746268656
$ws.Cells.Item(16, 13).Value = -37.857142857142
$ws.Cells.Item(16, 14).Value = -89.530685920577
$ws.Cells.Item(17, 3).Value = 9
$ws.Cells.Item(17, 4).Value = 18
$ws.Cells.Item(17, 5).Value = -50
$ws.Cells.Item(17, 6).Value = 49
$ws.Cells.Item(17, 7).Value = 45
$ws.Cells.Item(17, 8).Value = 8.888888888888
$ws.Cells.Item(17, 9).Value = 229
$ws.Cells.Item(17, 10).Value = 208
$ws.Cells.Item(17, 11).Value = 10.096153846153
$ws.Cells.Item(17, 12).Value = 16.243654822335
$ws.Cells.Item(17, 13).Value = 19.895287958115
$ws.Cells.Item(17, 14).Value = -35.854341736694
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 7
$ws.Cells.Item(18, 5).Value = -57.142857142857
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).Value = -38.888888888888
$ws.Cells.Item(18, 9).Value = 54
$ws.Cells.Item(18, 10).Value = 70
$ws.Cells.Item(18, 11).Value = -22.857142857142
$ws.Cells.Item(18, 12).Value = 35
$ws.Cells.Item(18, 13).Value = -45.454545454545
$ws.Cells.Item(18, 14).Value = -82.178217821782
$ws.Cells.Item(19, 3).Value = 9
$ws.Cells.Item(19, 4).Value = 7
$ws.Cells.Item(19, 5).Value = 28.571428571428
$ws.Cells.Item(19, 6).Value = 25
$ws.Cells.Item(19, 7).Value = 37
$ws.Cells.Item(19, 8).Value = -32.432432432432
$ws.Cells.Item(19, 9).Value = 119
$ws.Cells.Item(19, 10).Value = 150
$ws.Cells.Item(19, 11).Value = -20.666666666666
$ws.Cells.Item(19, 12).Value = 33.707865168539
$ws.Cells.Item(19, 13).Value = -12.5
$ws.Cells.Item(19, 14).Value = -58.680555555555
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = 50
$ws.Cells.Item(20, 7).Value = 8
$ws.Cells.Item(20, 8).Value = 12.5
$ws.Cells.Item(20, 9).Value = 40
$ws.Cells.Item(20, 10).Value = 77
$ws.Cells.Item(20, 11).Value = -48.051948051948
$ws.Cells.Item(20, 12).Value = 5.263157894736
$ws.Cells.Item(20, 13).Value = -29.824561403508
$ws.Cells.Item(20, 14).Value = -80.95238095238
$ws.Cells.Item(21, 3).Value = 33
$ws.Cells.Item(21, 4).Value = 40
$ws.Cells.Item(21, 5).Value = -17.5
$ws.Cells.Item(21, 6).Value = 114
$ws.Cells.Item(21, 7).Value = 139
$ws.Cells.Item(21, 8).Value = -17.985611510791
$ws.Cells.Item(21, 9).Value = 544
$ws.Cells.Item(21, 10).Value = 622
$ws.Cells.Item(21, 11).Value = -12.540192926045
$ws.Cells.Item(21, 12).Value = 19.298245614035
$ws.Cells.Item(21, 13).Value = -15.132605304212
$ws.Cells.Item(21, 14).Value = -73.359451518119
$ws.Cells.Item(22, 6).Value = 3
$ws.Cells.Item(22, 7).Value = 7
$ws.Cells.Item(22, 8).Value = -57.142857142857
$ws.Cells.Item(22, 9).Value = 18
$ws.Cells.Item(22, 10).Value = 24
$ws.Cells.Item(22, 11).Value = -25
$ws.Cells.Item(22, 12).Value = 20
$ws.Cells.Item(22, 13).Value = 12.5
$ws.Cells.Item(23, 3).Value = 5
$ws.Cells.Item(23, 4).Value = 7
$ws.Cells.Item(23, 5).Value = -28.571428571428
$ws.Cells.Item(23, 6).Value = 22
$ws.Cells.Item(23, 7).Value = 24
$ws.Cells.Item(23, 8).Value = -8.333333333333
$ws.Cells.Item(23, 9).Value = 117
$ws.Cells.Item(23, 10).Value = 123
$ws.Cells.Item(23, 11).Value = -4.878048780487
$ws.Cells.Item(23, 12).Value = 11.428571428571
$ws.Cells.Item(23, 13).Value = 64.788732394366
$ws.Cells.Item(24, 3).Value = 20
$ws.Cells.Item(24, 4).Value = 25
$ws.Cells.Item(24, 5).Value = -20
$ws.Cells.Item(24, 6).Value = 72
$ws.Cells.Item(24, 7).Value = 78
$ws.Cells.Item(24, 8).Value = -7.692307692307
$ws.Cells.Item(24, 9).Value = 423
$ws.Cells.Item(24, 10).Value = 372
$ws.Cells.Item(24, 11).Value = 13.709677419354
$ws.Cells.Item(24, 12).Value = 56.666666666666
$ws.Cells.Item(24, 13).Value = 38.235294117647
$ws.Cells.Item(25, 3).Value = 23
$ws.Cells.Item(25, 4).Value = 22
$ws.Cells.Item(25, 5).Value = 4.545454545454
$ws.Cells.Item(25, 6).Value = 69
$ws.Cells.Item(25, 7).Value = 97
$ws.Cells.Item(25, 8).Value = -28.865979381443
$ws.Cells.Item(25, 9).Value = 326
$ws.Cells.Item(25, 10).Value = 347
$ws.Cells.Item(25, 11).Value = -6.051873198847
$ws.Cells.Item(25, 12).Value = 17.266187050359
$ws.Cells.Item(25, 13).Value = -28.508771929824
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = -100
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = -66.666666666666
$ws.Cells.Item(26, 10).Value = 14
$ws.Cells.Item(26, 11).Value = 14.285714285714
$ws.Cells.Item(26, 12).Value = -40.74074074074
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = -16.666666666666
$ws.Cells.Item(27, 10).Value = 23
$ws.Cells.Item(27, 11).Value = 34.782608695652
$ws.Cells.Item(27, 12).Value = -3.125
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 8).Value = -75
$ws.Cells.Item(28, 10).Value = 21
$ws.Cells.Item(28, 11).Value = -38.095238095238
$ws.Cells.Item(28, 12).Value = -53.571428571428
$ws.Cells.Item(28, 13).Value = -53.571428571428
$ws.Cells.Item(28, 14).Value = -87.128712871287
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 8).Value = -66.666666666666
$ws.Cells.Item(29, 10).Value = 18
$ws.Cells.Item(29, 11).Value = -33.333333333333
$ws.Cells.Item(29, 12).Value = -45.454545454545
$ws.Cells.Item(29, 13).Value = -40
$ws.Cells.Item(29, 14).Value = -87.5
